$wb = $excel.ActiveWorkbook

# --- Worksheets ---
$wsHS = $wb.Worksheets.Item(1)    # "HS-422"
$wsM1 = $wb.Worksheets.Item(2)    # "M1" -> will become "M0"

# --- M1 sheet data: new linear-equation sample values ---
# Row 1 stays x=0, y changes 625000 -> 509167
$wsM1.Range("B1").Value = 509167
# Row 2 becomes x=180 (was 90), y=2304027 (was 1450000)
$wsM1.Range("A2").Value = 180
$wsM1.Range("B2").Value = 2304027
# Old row 3 (x=180, y=2250000) is removed entirely
$wsM1.Range("A3:B3").ClearContents()

# --- New styled (centered + wrapped) blank cells A6:C6 on the M1 sheet ---
$a6 = $wsM1.Range("A6")
$a6.HorizontalAlignment = -4108   # xlCenter
$a6.WrapText = $true
$a6.Copy()
$wsM1.Range("B6:C6").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# --- Move the M1 chart up by one row (default row height = 16pt) ---
$co = $wsM1.ChartObjects().Item(1)
$co.Top = $co.Top - 16

# --- Sheet view changes ---
# HS-422: zoom out + pin top-left cell
$wsHS.Application.ActiveWindow
$wsHS.Activate()
$excel.ActiveWindow.Zoom = 64
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 2
$wsHS.Range("N9").Select()

# M1: zoom out + change selection to B11
$wsM1.Activate()
$excel.ActiveWindow.Zoom = 160
$wsM1.Range("B11").Select()

# --- Rename sheet M1 -> M0 (after all M1-keyed lookups above) ---
$wsM1.Name = "M0"

# --- Workbook window size ---
$excel.ActiveWindow.Width = 18220
$excel.ActiveWindow.Height = 16720

# Re-activate M0 (tab 2) to match the saved activeTab/tabSelected state
$wsM1.Activate()
